{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2024-03-23 Saturday\", \"2024-03-24 Sunday\"],\n  [\"46\u00f73=15, 1\", \"30\u00f78=3, 6\"],\n  [\"68\u00f77=9, 5\", \"88\u00f75=17, 3\"],\n  [\"86\u00f78=10, 6\", \"72\u00f73=24, 0\"],\n  [\"43\u00f78=5, 3\", \"32\u00f77=4, 4\"],\n  [\"19\u00f75=3, 4\", \"64\u00f75=12, 4\"],\n  [\"17\u00f76=2, 5\", \"12\u00f77=1, 5\"],\n  [\"21\u00f76=3, 3\", \"54\u00f74=13, 2\"],\n  [\"67\u00f79=7, 4\", \"32\u00f73=10, 2\"],\n  [\"88\u00f72=44, 0\", \"17\u00f79=1, 8\"],\n  [\"98\u00f74=24, 2\", \"19\u00f76=3, 1\"],\n  [\"26\u00f74=6, 2\", \"31\u00f78=3, 7\"],\n  [\"80\u00f78=10, 0\", \"14\u00f73=4, 2\"],\n  [\"10\u00f73=3, 1\", \"77\u00f75=15, 2\"],\n  [\"30\u00f74=7, 2\", \"91\u00f76=15, 1\"],\n  [\"23\u00f78=2, 7\", \"66\u00f74=16, 2\"],\n  [\"39\u00f77=5, 4\", \"47\u00f72=23, 1\"],\n  [\"13\u00f79=1, 4\", \"24\u00f76=4, 0\"],\n  [\"14\u00f75=2, 4\", \"78\u00f74=19, 2\"],\n  [\"81\u00f73=27, 0\", \"96\u00f78=12, 0\"],\n  [\"46\u00f74=11, 2\", \"72\u00f74=18, 0\"],\n  [\"14\u00f78=1, 6\", \"84\u00f79=9, 3\"],\n  [\"47\u00f79=5, 2\", \"12\u00f73=4, 0\"],\n  [\"61\u00f72=30, 1\", \"20\u00f79=2, 2\"],\n  [\"63\u00f75=12, 3\", \"94\u00f77=13, 3\"],\n  [\"21\u00f77=3, 0\", \"19\u00f74=4, 3\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$null = $d.Content.Find.Execute(\"2024-03-23 Saturday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-03-24 Sunday\", 2)\n$null = $d.Content.Find.Execute(\"46\u00f73=15, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"30\u00f78=3, 6\", 2)\n$null = $d.Content.Find.Execute(\"68\u00f77=9, 5\", $false, $false, $false, $false, $false, $true, 1, $false, \"88\u00f75=17, 3\", 2)\n$null = $d.Content.Find.Execute(\"86\u00f78=10, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"72\u00f73=24, 0\", 2)\n$null = $d.Content.Find.Execute(\"43\u00f78=5, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"32\u00f77=4, 4\", 2)\n$null = $d.Content.Find.Execute(\"19\u00f75=3, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"64\u00f75=12, 4\", 2)\n$null = $d.Content.Find.Execute(\"17\u00f76=2, 5\", $false, $false, $false, $false, $false, $true, 1, $false, \"12\u00f77=1, 5\", 2)\n$null = $d.Content.Find.Execute(\"21\u00f76=3, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"54\u00f74=13, 2\", 2)\n$null = $d.Content.Find.Execute(\"67\u00f79=7, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"32\u00f73=10, 2\", 2)\n$null = $d.Content.Find.Execute(\"88\u00f72=44, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"17\u00f79=1, 8\", 2)\n$null = $d.Content.Find.Execute(\"98\u00f74=24, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00f76=3, 1\", 2)\n$null = $d.Content.Find.Execute(\"26\u00f74=6, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"31\u00f78=3, 7\", 2)\n$null = $d.Content.Find.Execute(\"80\u00f78=10, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"14\u00f73=4, 2\", 2)\n$null = $d.Content.Find.Execute(\"10\u00f73=3, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"77\u00f75=15, 2\", 2)\n$null = $d.Content.Find.Execute(\"30\u00f74=7, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"91\u00f76=15, 1\", 2)\n$null = $d.Content.Find.Execute(\"23\u00f78=2, 7\", $false, $false, $false, $false, $false, $true, 1, $false, \"66\u00f74=16, 2\", 2)\n$null = $d.Content.Find.Execute(\"39\u00f77=5, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"47\u00f72=23, 1\", 2)\n$null = $d.Content.Find.Execute(\"13\u00f79=1, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"24\u00f76=4, 0\", 2)\n$null = $d.Content.Find.Execute(\"14\u00f75=2, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"78\u00f74=19, 2\", 2)\n$null = $d.Content.Find.Execute(\"81\u00f73=27, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"96\u00f78=12, 0\", 2)\n$null = $d.Content.Find.Execute(\"46\u00f74=11, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"72\u00f74=18, 0\", 2)\n$null = $d.Content.Find.Execute(\"14\u00f78=1, 6\", $false, $false, $false, $false, $false, $true, 1, $false, \"84\u00f79=9, 3\", 2)\n$null = $d.Content.Find.Execute(\"47\u00f79=5, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"12\u00f73=4, 0\", 2)\n$null = $d.Content.Find.Execute(\"61\u00f72=30, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"20\u00f79=2, 2\", 2)\n$null = $d.Content.Find.Execute(\"63\u00f75=12, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"94\u00f77=13, 3\", 2)\n$null = $d.Content.Find.Execute(\"21\u00f77=3, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00f74=4, 3\", 2)\n"}
